$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.902.82'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.632.78'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.73'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.69'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0903'
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.868.83'
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").Value = '1.640.23'
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.28'
$ws.Range("E15").Value = '  +6.79%  '
$ws.Range("D16").Value = '29.935.52'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.84'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.10'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.29'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '0.0₃0702'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.76'
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  +2.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.09'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.45'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.57'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("E31").Value = '  +3.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '1.424.53'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("E35").Value = '  +4.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -1.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.76'
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.50'
$ws.Range("E40").Value = '  +10.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.553'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0500'
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.98'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.827'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '50.73'
$ws.Range("E47").Value = '  -6.43%  '
$ws.Range("D48").Value = '1.776.36'
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.34'
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.63'
$ws.Range("E50").Value = '  +3.84%  '
$ws.Range("E51").Value = '  +10.39%  '
